# Update "想去人数" (interest count) figures in both the "展览" and
# "全部类型" sheets to match the refreshed scrape output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 382
$ws1.Range("F4").Value = 4887
$ws1.Range("F5").Value = 23
$ws1.Range("F6").Value = 22
$ws1.Range("F8").Value = 487

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 382
$ws4.Range("F4").Value = 4887
$ws4.Range("F6").Value = 23
$ws4.Range("F7").Value = 22
$ws4.Range("F10").Value = 487
